# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
# Swap the data (columns B:AC) between pairs of rows that got reordered
# in the source feed. Column A (the running index) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(17, 18),
    @(22, 23),
    @(69, 70),
    @(78, 79),
    @(91, 92),
    @(103, 104),
    @(108, 109)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B${r1}:AC${r1}")
    $rng2 = $ws.Range("B${r2}:AC${r2}")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
